# ---------------------------------------------------------------------------
# 19/09/2017 MAMATHA CHITRA CHICK IN
#
# 1) The "Mon Sep 17 11:43:26 PDT 2017" timestamp was previously split across
#    two runs; collapse it back into a single run (same visible text/format).
# 2) Append a brand-new "Tue Sep 18 11:26:16 PDT 2017" purchase-details block
#    (SHANTHARAJU / CHOW / CARROT EVE) after the last entry in the document.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- Change 1: merge the two date/time runs into a single run -------------
$d.Content.Find.Execute(
    "Mon Sep 17 11:43:26 PDT 2017", $true, $false, $false, $false, $false,
    $true, 1, $false, "Mon Sep 17 11:43:26 PDT 2017", 2) | Out-Null

# --- Change 2: build the new "Tue Sep 18" block ----------------------------

$wNs  = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$rPr  = '<w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr>'
$rPrB = '<w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:b/></w:rPr>'
$rPrR = '<w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="FF0000"/></w:rPr>'

function New-EmptyParagraph([bool]$bold = $false) {
    $pPrRpr = if ($bold) { $rPrB } else { $rPr }
    return "<w:p><w:pPr><w:pStyle w:val=`"PlainText`"/>$pPrRpr</w:pPr></w:p>"
}

# A "label <tabs...> - value" line, matching the document's existing layout.
function New-FieldParagraph([string]$label, [int]$tabCount, [string]$value, [string]$runProps = $rPr) {
    $p = "<w:p><w:pPr><w:pStyle w:val=`"PlainText`"/>$runProps</w:pPr>"
    $p += "<w:r>$runProps<w:t>$label</w:t></w:r>"
    for ($i = 1; $i -lt $tabCount; $i++) {
        $p += "<w:r>$runProps<w:tab/></w:r>"
    }
    $p += "<w:r>$runProps<w:tab/><w:t>$value</w:t></w:r>"
    $p += "</w:p>"
    return $p
}

function New-TextParagraph([string]$text) {
    return "<w:p><w:pPr><w:pStyle w:val=`"PlainText`"/>$rPr</w:pPr><w:r>$rPr<w:t>$text</w:t></w:r></w:p>"
}

$dateParagraph = "<w:p><w:pPr><w:pStyle w:val=`"PlainText`"/>$rPr</w:pPr>" +
    "<w:r>$rPr<w:t>Tue Sep 18</w:t></w:r>" +
    "<w:r>$rPr<w:t xml:space=`"preserve`"> 11:26:16 PDT 2017</w:t></w:r></w:p>"

$blockParts = @(
    (New-EmptyParagraph($true)),
    $dateParagraph,
    (New-FieldParagraph "Person Name" 4 "- SHANTHARAJU"),
    (New-TextParagraph "---------------------------------------------------------------"),
    (New-FieldParagraph "Item Name" 4 "- CHOW"),
    (New-FieldParagraph "Amount Received" 3 "- 3660" $rPrR),
    (New-FieldParagraph "Amount Received mode" 2 "- CASH AND CLEAR"),
    (New-EmptyParagraph($false)),
    (New-FieldParagraph "Item Name" 4 "- CARROT EVE"),
    (New-FieldParagraph "Number of Pockets" 3 "- 2"),
    (New-FieldParagraph "Number of KGs" 3 "- 165"),
    (New-FieldParagraph "Rate" 5 "- 20"),
    (New-FieldParagraph "Total Price" 4 "- 3300.0"),
    (New-FieldParagraph "Amount balance" 3 "- 3300.0" $rPrB),
    (New-EmptyParagraph($false)),
    (New-EmptyParagraph($true))
)

$newXml = [string]::Join("", $blockParts)
# The fragment needs its own namespace context the first time it is parsed.
$newXml = $newXml -replace "^<w:p>", "<w:p $wNs>"

# --- Locate the last "Amount balance" paragraph (the final entry) and -----
# --- insert the new block immediately after it. ----------------------------
$rng = $d.Content
$lastMatch = $null
while ($rng.Find.Execute("Amount balance", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $lastMatch = $d.Range($rng.Start, $rng.End)
    $rng.Collapse(0)
}

$para = $lastMatch.Paragraphs.Last
$insPoint = $d.Range($para.Range.End, $para.Range.End)
[void]$insPoint.InsertXML($newXml)
